$wb = $excel.ActiveWorkbook

# Sheet ALC, row 46
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 1000
$ws.Range("J46").Value = 1000
$ws.Range("L46").Value = 3000
$ws.Range("N46").Value = -3238

# Sheet ALC, row 60
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H60").Value = 1000
$ws.Range("J60").Value = 1000
$ws.Range("L60").Value = 3000
$ws.Range("N60").Value = -3968

# Sheet ALC, row 121
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 2005
$ws.Range("J121").Value = 2005
$ws.Range("L121").Value = 6015
$ws.Range("N121").Value = -9509

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1330.8889
$ws.Range("I137").Value = 1327.091
$ws.Range("J137").Value = 1336.8572
$ws.Range("K137").Value = 3981.273
$ws.Range("L137").Value = 4010.5716
$ws.Range("M137").Value = -1431.273
$ws.Range("N137").Value = -9110.571599999999

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9055.267
$ws.Range("I32").Value = 9523.5
$ws.Range("J32").Value = 2500
$ws.Range("K32").Value = 9523.5
$ws.Range("L32").Value = 2500
$ws.Range("M32").Value = -9236.5
$ws.Range("N32").Value = -3074

# Sheet ARM, row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1611.6923
$ws.Range("I45").Value = 1725.3
$ws.Range("J45").Value = 1233
$ws.Range("K45").Value = 1725.3
$ws.Range("L45").Value = 1233
$ws.Range("M45").Value = -1348.3
$ws.Range("N45").Value = -1987

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1217.129
$ws.Range("I61").Value = 873.7308
$ws.Range("K61").Value = 873.7308
$ws.Range("M61").Value = -661.7308

# Sheet ARM, row 75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

# Sheet ARM, row 78
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2117.88
$ws.Range("I132").Value = 1697.5
$ws.Range("J132").Value = 3799.4
$ws.Range("K132").Value = 5092.5
$ws.Range("L132").Value = 11398.2
$ws.Range("M132").Value = -2562.5
$ws.Range("N132").Value = -16458.2

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1217.129
$ws.Range("I136").Value = 873.7308
$ws.Range("K136").Value = 2621.1924
$ws.Range("M136").Value = -71.19239999999991

# Sheet BSM, row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1861.3077
$ws.Range("I20").Value = 1838.5
$ws.Range("K20").Value = 1838.5
$ws.Range("M20").Value = -1591.5

# Sheet BSM, row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4378.4116
$ws.Range("I86").Value = 4763.5386
$ws.Range("J86").Value = 3126.75
$ws.Range("K86").Value = 4763.5386
$ws.Range("L86").Value = 3126.75
$ws.Range("M86").Value = -3640.5386
$ws.Range("N86").Value = -5372.75

# Sheet BSM, row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4378.4116
$ws.Range("I89").Value = 4763.5386
$ws.Range("J89").Value = 3126.75
$ws.Range("K89").Value = 23817.693
$ws.Range("L89").Value = 15633.75
$ws.Range("M89").Value = -18201.693
$ws.Range("N89").Value = -26865.75

# Sheet BSM, row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 22728250
$ws.Range("I94").Value = 35715196
$ws.Range("J94").Value = 1095
$ws.Range("K94").Value = 35715196
$ws.Range("L94").Value = 1095
$ws.Range("M94").Value = -35714745
$ws.Range("N94").Value = -1997

# Sheet BSM, row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2019.8
$ws.Range("I107").Value = 1499.5
$ws.Range("J107").Value = 2366.6667
$ws.Range("K107").Value = 1499.5
$ws.Range("L107").Value = 2366.6667
$ws.Range("M107").Value = 420.5
$ws.Range("N107").Value = -6206.6667

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1712.9333
$ws.Range("I31").Value = 855.14813
$ws.Range("K31").Value = 855.14813
$ws.Range("M31").Value = -560.14813

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1712.9333
$ws.Range("I34").Value = 855.14813
$ws.Range("K34").Value = 855.14813
$ws.Range("M34").Value = -653.14813

# Sheet CRP, row 44
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

# Sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1058.826
$ws.Range("I58").Value = 869.64703
$ws.Range("J58").Value = 1594.8334
$ws.Range("K58").Value = 869.64703
$ws.Range("L58").Value = 1594.8334
$ws.Range("M58").Value = -666.64703
$ws.Range("N58").Value = -2000.8334

# Sheet CRP, row 80
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 15000
$ws.Range("J80").Value = 15000
$ws.Range("L80").Value = 15000
$ws.Range("N80").Value = -17246

# Sheet CRP, row 83
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H83").Value = 15000
$ws.Range("J83").Value = 15000
$ws.Range("L83").Value = 45000
$ws.Range("N83").Value = -56232

# Sheet CRP, row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 8361836
$ws.Range("I86").Value = 16670068
$ws.Range("K86").Value = 16670068
$ws.Range("M86").Value = -16668945

# Sheet CRP, row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 8361836
$ws.Range("I89").Value = 16670068
$ws.Range("K89").Value = 83350340
$ws.Range("M89").Value = -83344724

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2970.2144
$ws.Range("I132").Value = 2358.5
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 7075.5
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -4545.5
$ws.Range("N132").Value = -18558.5

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 15874057
$ws.Range("I134").Value = 17544878
$ws.Range("J134").Value = 1257
$ws.Range("K134").Value = 52634634
$ws.Range("L134").Value = 3771
$ws.Range("M134").Value = -52632099
$ws.Range("N134").Value = -8841

# Sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1058.826
$ws.Range("I136").Value = 869.64703
$ws.Range("J136").Value = 1594.8334
$ws.Range("K136").Value = 2608.94109
$ws.Range("L136").Value = 4784.5002
$ws.Range("M136").Value = -58.9410899999998
$ws.Range("N136").Value = -9884.5002

# Sheet CUL, row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3300
$ws.Range("J55").Value = 3300
$ws.Range("L55").Value = 9900
$ws.Range("N55").Value = -10254

# Sheet CUL, row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 809.6
$ws.Range("I107").Value = 150
$ws.Range("J107").Value = 911.0769
$ws.Range("K107").Value = 450
$ws.Range("L107").Value = 2733.2307
$ws.Range("M107").Value = 1470
$ws.Range("N107").Value = -6573.2307

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 28573010
$ws.Range("J131").Value = 1671.5454
$ws.Range("L131").Value = 5014.6362
$ws.Range("N131").Value = -15094.6362

# Sheet CUL, row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1750
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 13500
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -10970
$ws.Range("N132").Value = -23060

# Sheet GSM, row 14
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 7000501
$ws.Range("I14").Value = 9333335
$ws.Range("J14").Value = 2000
$ws.Range("K14").Value = 9333335
$ws.Range("L14").Value = 2000
$ws.Range("M14").Value = -9333167
$ws.Range("N14").Value = -2336

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2035.3334
$ws.Range("I132").Value = 1514.9412
$ws.Range("J132").Value = 4247
$ws.Range("K132").Value = 4544.8236
$ws.Range("L132").Value = 12741
$ws.Range("M132").Value = -2014.8236
$ws.Range("N132").Value = -17801

# Sheet GSM, row 134
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 16753.857
$ws.Range("J134").Value = 16753.857
$ws.Range("L134").Value = 50261.571
$ws.Range("N134").Value = -55331.571

# Sheet LTW, row 31
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 3573.889
$ws.Range("I31").Value = 1115
$ws.Range("J31").Value = 4276.4287
$ws.Range("K31").Value = 1115
$ws.Range("L31").Value = 4276.4287
$ws.Range("M31").Value = -867
$ws.Range("N31").Value = -4772.4287

# Sheet LTW, row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 11500
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 20000
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 20000
$ws.Range("M40").Value = -2864
$ws.Range("N40").Value = -20272

# Sheet LTW, row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2215.3845
$ws.Range("I100").Value = 1890
$ws.Range("K100").Value = 1890
$ws.Range("M100").Value = -1349

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 19233790
$ws.Range("I122").Value = 27780556
$ws.Range("J122").Value = 3566
$ws.Range("K122").Value = 83341668
$ws.Range("L122").Value = 10698
$ws.Range("M122").Value = -83339218
$ws.Range("N122").Value = -15598

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 31762.242
$ws.Range("I132").Value = 993.5833
$ws.Range("K132").Value = 2980.7499
$ws.Range("M132").Value = -450.7498999999998

# Sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2212.75
$ws.Range("I136").Value = 1926
$ws.Range("K136").Value = 5778
$ws.Range("M136").Value = -3228

# Sheet WVR, row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 559.2
$ws.Range("I107").Value = 449
$ws.Range("K107").Value = 1347
$ws.Range("M107").Value = 573

# Sheet WVR, row 119
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 43348
$ws.Range("J119").Value = 43348
$ws.Range("L119").Value = 43348
$ws.Range("N119").Value = -53024

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 54785024
$ws.Range("I122").Value = 57275210
$ws.Range("K122").Value = 171825630
$ws.Range("M122").Value = -171823180

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 747.4074000000001
$ws.Range("I136").Value = 676.0476
$ws.Range("J136").Value = 997.1667
$ws.Range("K136").Value = 2028.1428
$ws.Range("L136").Value = 2991.5001
$ws.Range("M136").Value = 521.8571999999999
$ws.Range("N136").Value = -8091.5001
